# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-33 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 33 (in row order)
$newValues = @(1, 0, 0, 0, 0, 2, 0, 0, 2, 0, 2, 2, 0, 0, 0, 0, 0, 1, 3, 2, 2, 2, 1, 1, 2, 0, 0, 0, 0, 2, 1, 0)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
